$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$formula = "=(TestingFile[[#This Row],[Obj. LR]]-TestingFile[[#This Row],[LB Heuristic]])/TestingFile[[#This Row],[Obj. LR]]"

for ($r = 2; $r -le 121; $r++) {
    $ws.Range("G$r").Formula = $formula
}

$ws.Range("G124").Select()

Write-Output "done"
